$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.Text = "Header text"
"done"
